# Update "想去人数" (want-to-go count) values in column F across the four
# sheets of the workbook, per the source diff.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 (Exhibition) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value  = 626
$ws.Range("F5").Value  = 2775
$ws.Range("F9").Value  = 6357
$ws.Range("F13").Value = 5025
$ws.Range("F14").Value = 105
$ws.Range("F16").Value = 2626
$ws.Range("F17").Value = 1352
$ws.Range("F18").Value = 1517
$ws.Range("F20").Value = 315
$ws.Range("F21").Value = 121
$ws.Range("F23").Value = 1069
$ws.Range("F24").Value = 241
$ws.Range("F26").Value = 535
$ws.Range("F27").Value = 1378
$ws.Range("F28").Value = 1047
$ws.Range("F31").Value = 580
$ws.Range("F32").Value = 29
$ws.Range("F33").Value = 27
$ws.Range("F34").Value = 90
$ws.Range("F35").Value = 247
$ws.Range("F36").Value = 1497
$ws.Range("F39").Value = 115
$ws.Range("F41").Value = 18
$ws.Range("F42").Value = 298
$ws.Range("F43").Value = 2283
$ws.Range("F44").Value = 2555
$ws.Range("F48").Value = 103

# --- Sheet 2: 演出 (Performance) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F20").Value = 37
$ws.Range("F24").Value = 378
$ws.Range("F25").Value = 30

# --- Sheet 3: 本地生活 (Local Life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F6").Value  = 1700
$ws.Range("F8").Value  = 1501
$ws.Range("F10").Value = 2526
$ws.Range("F11").Value = 854
$ws.Range("F13").Value = 29

# --- Sheet 4: 全部类型 (All Types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value  = 626
$ws.Range("F7").Value  = 2775
$ws.Range("F9").Value  = 1501
$ws.Range("F11").Value = 2526
$ws.Range("F12").Value = 6357
$ws.Range("F13").Value = 854
$ws.Range("F16").Value = 5025
$ws.Range("F17").Value = 2626
$ws.Range("F18").Value = 1352
$ws.Range("F19").Value = 1517
$ws.Range("F21").Value = 121
$ws.Range("F24").Value = 241
$ws.Range("F27").Value = 1378
$ws.Range("F28").Value = 1047
$ws.Range("F31").Value = 580
$ws.Range("F32").Value = 29
$ws.Range("F34").Value = 27
$ws.Range("F35").Value = 247
$ws.Range("F39").Value = 18
$ws.Range("F42").Value = 298
$ws.Range("F43").Value = 30
$ws.Range("F44").Value = 2283
$ws.Range("F45").Value = 2555
$ws.Range("F48").Value = 103
